# The "client" sheet had a row that was imported with every column set to
# null/blank. Excel represents this as a new row whose only populated cell
# is a single space character in column B (the rest stay empty).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("client")

$ws.Range("B3").Value = " "
$ws.Range("B3").Select()
